$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original (default) style/number-format of the Price column,
# switch it to Text temporarily so that numeric-looking strings (e.g. "0.700",
# "173.20") are written verbatim instead of being parsed into floating point
# numbers, then restore the original style afterwards.
$origPriceStyle = $ws.Range("D2").Style
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.805.56'
$ws.Range('E2').Value = '  -7.43%  '

$ws.Range('D3').Value = '3.680.84'
$ws.Range('E3').Value = '  -7.57%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '568.81'
$ws.Range('E5').Value = '  -6.65%  '

$ws.Range('D6').Value = '173.20'
$ws.Range('E6').Value = '  +0.41%  '

$ws.Range('D7').Value = '3.667.00'
$ws.Range('E7').Value = '  -7.71%  '

$ws.Range('D8').Value = '0.622'
$ws.Range('E8').Value = '  -11.48%  '

$ws.Range('E9').Value = '  +0.15%  '

$ws.Range('D10').Value = '0.700'
$ws.Range('E10').Value = '  -12.76%  '

$ws.Range('E11').Value = '  -13.23%  '

$ws.Range('D12').Value = '51.16'
$ws.Range('E12').Value = '  -10.31%  '

$ws.Range('E13').Value = '  -13.48%  '

$ws.Range('D14').Value = '10.43'
$ws.Range('E14').Value = '  -11.75%  '

$ws.Range('D15').Value = '4.255.06'
$ws.Range('E15').Value = '  -7.80%  '

$ws.Range('D16').Value = '3.674.85'
$ws.Range('E16').Value = '  -7.81%  '

$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '0.126'
$ws.Range('E17').Value = '  -3.42%  '

$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '19.28'
$ws.Range('E18').Value = '  -8.75%  '

$ws.Range('D19').Value = '12.78'
$ws.Range('E19').Value = '  -11.00%  '

$ws.Range('E20').Value = '  -11.12%  '

$ws.Range('D21').Value = '67.465.84'
$ws.Range('E21').Value = '  -7.74%  '

$ws.Range('D22').Value = '404.43'
$ws.Range('E22').Value = '  -13.80%  '

$ws.Range('D23').Value = '4.42'
$ws.Range('E23').Value = '  -8.59%  '

$ws.Range('D24').Value = '87.46'
$ws.Range('E24').Value = '  -10.51%  '

$ws.Range('D25').Value = '3.03'
$ws.Range('E25').Value = '  -11.47%  '

$ws.Range('D26').Value = '12.69'
$ws.Range('E26').Value = '  -11.68%  '

$ws.Range('D27').Value = '10.60'
$ws.Range('E27').Value = '  -6.36%  '

$ws.Range('D28').Value = '5.97'
$ws.Range('E28').Value = '  +1.26%  '

$ws.Range('D29').Value = '3.71'
$ws.Range('E29').Value = '  -12.51%  '

$ws.Range('D30').Value = '9.40'
$ws.Range('E30').Value = '  -12.25%  '

$ws.Range('D31').Value = '32.49'
$ws.Range('E31').Value = '  -11.34%  '

$ws.Range('D32').Value = '7.57'
$ws.Range('E32').Value = '  -6.13%  '

$ws.Range('D33').Value = '12.44'
$ws.Range('E33').Value = '  -12.15%  '

$ws.Range('D34').Value = '0.116'
$ws.Range('E34').Value = '  -11.34%  '

$ws.Range('D35').Value = '64.54'
$ws.Range('E35').Value = '  -8.70%  '

$ws.Range('D36').Value = '42.81'
$ws.Range('E36').Value = '  -14.81%  '

$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').Value = '586.73'
$ws.Range('E37').Value = '  -8.75%  '

$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0888'
$ws.Range('E38').Value = '  -13.72%  '

$ws.Range('E39').Value = '  +0.04%  '

$ws.Range('D40').Value = '0.394'
$ws.Range('E40').Value = '  -9.93%  '

$ws.Range('E41').Value = '  -0.06%  '

$ws.Range('E42').Value = '  -10.51%  '

$ws.Range('D43').Value = '2.98'
$ws.Range('E43').Value = '  -9.81%  '

$ws.Range('D44').Value = '2.96'
$ws.Range('E44').Value = '  -13.40%  '

$ws.Range('E45').Value = '  -11.03%  '

$ws.Range('D46').Value = '2.57'
$ws.Range('E46').Value = '  -3.57%  '

$ws.Range('D47').Value = '9.14'
$ws.Range('E47').Value = '  -13.75%  '

$ws.Range('E48').Value = '  -12.41%  '

$ws.Range('D49').Value = '2.69'
$ws.Range('E49').Value = '  -10.79%  '

$ws.Range('D50').Value = '3.12'
$ws.Range('E50').Value = '  -9.47%  '

$ws.Range('D51').Value = '2.718.50'
$ws.Range('E51').Value = '  -3.84%  '

# Restore the original style on the Price column now that the text values
# have been written, so cell formatting matches the source workbook.
$ws.Range("D2:D51").Style = $origPriceStyle
